# "Fixing wrongly done merge"
# week1 sheet had a duplicated/misplaced "Meeting semanal" entry (row 9) that was
# really meant to read "Meetings", and the last task row (13) had lost its text
# during a bad merge. Restore the correct contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")

# Row 9: rename task text from "Meeting semanal" to "Meetings"
$ws.Range("B9").Value = "Meetings"

# Row 13: restore the missing task text (previously empty, underline style)
# and give it the normal (non-underlined) style used by the other task rows,
# matching cell F9's style.
$ws.Range("B13").Value = "Fazer pdf dos 3 User Stories mais votados pela equipa e submeter no moodle"
$ws.Range("B13").Style = $ws.Range("F9").Style
